# Fruta / hortaliza, semanal
# Insert a new weekly record as row 85 (Terminal Hortofrutícola Agro Chillán - Perejil),
# pushing the existing rows 85..142 down to 86..143.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 85 (shifts rows 85-142 down to 86-143).
$ws.Rows.Item(85).Insert()

# Populate the new row 85 with the new weekly observation.
$ws.Cells.Item(85, 1).Value = 7
$ws.Cells.Item(85, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(85, 3).Value = "Ñuble"
$ws.Cells.Item(85, 4).Value = 45236
$ws.Cells.Item(85, 5).Value = 16
$ws.Cells.Item(85, 6).Value = 100112044
$ws.Cells.Item(85, 7).Value = "Perejil"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 120
$ws.Cells.Item(85, 11).Value = 2000
$ws.Cells.Item(85, 12).Value = 2000
$ws.Cells.Item(85, 13).Value = 2000
$ws.Cells.Item(85, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(85, 15).Value = "Región de Ñuble"
$ws.Cells.Item(85, 16).Value = 2000
$ws.Cells.Item(85, 17).Value = 1
$ws.Cells.Item(85, 18).Value = "Hortaliza"

# Match the style (date number format) used by column D in the surrounding rows.
$ws.Cells.Item(85, 4).NumberFormat = $ws.Cells.Item(86, 4).NumberFormat
